$p = $ppt.ActivePresentation

# --- Slide 13: "M16: Het project gebruikt tools voor vastgestelde taken" ---
# Remove the bullet paragraph about "controleren van de configuratie op
# aanwezigheid van bekende kwetsbaarheden," (OpenVAS reference removed).
$s13 = $p.Slides.Item(13)
$shape13 = $s13.Shapes.Item(2)
$tr13 = $shape13.TextFrame2.TextRange
$targetText = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,"
$paraCount13 = $tr13.Paragraphs().Count
for ($i = $paraCount13; $i -ge 1; $i--) {
    $para = $tr13.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd([char]13)
    if ($paraText -eq $targetText) {
        $para.Delete()
    }
}

# --- Slide 19: "M23: Het project zorgt voor de aanwezigheid van ... de Kwaliteitsaanpak" ---
$s19 = $p.Slides.Item(19)

# Title: add "kennis van en " before "ervaring met de Kwaliteitsaanpak"
$titleShape = $s19.Shapes.Item(1)
$titleRun = $titleShape.TextFrame2.TextRange.Runs(1, 1)
$titleRun.Text = "M23: Het project zorgt voor de aanwezigheid van kennis van en ervaring met de Kwaliteitsaanpak"

# Body: append explanation sentence about new project members
$bodyShape = $s19.Shapes.Item(2)
$bodyRun = $bodyShape.TextFrame2.TextRange.Runs(1, 1)
$bodyRun.Text = "De software delivery manager zorgt ervoor dat bij nieuwe projecten wordt gestart met ten minste twee projectleden die bekend zijn met de Kwaliteitsaanpak. Projectleden die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak."
